$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0

# Row 51
$ws.Range("H51").Value = 1639.8
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1639.8
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 1639.8
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2607.8

# Row 55
$ws.Range("H55").Value = 82.07692
$ws.Range("I55").Value = 79.09999999999999
$ws.Range("J55").Value = 92
$ws.Range("K55").Value = 79.09999999999999
$ws.Range("L55").Value = 92
$ws.Range("M55").Value = 134.9
$ws.Range("N55").Value = -520

# Row 74
$ws.Range("H74").Value = 4146.826
$ws.Range("I74").Value = 3522.2222
$ws.Range("K74").Value = 3522.2222
$ws.Range("M74").Value = -2586.2222

# Row 77
$ws.Range("H77").Value = 4146.826
$ws.Range("I77").Value = 3522.2222
$ws.Range("K77").Value = 17611.111
$ws.Range("M77").Value = -12931.111

# Row 113
$ws.Range("H113").Value = 2752.6667
$ws.Range("J113").Value = 2752.75
$ws.Range("L113").Value = 2752.75
$ws.Range("N113").Value = -9260.75

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Row 55
$ws.Range("H55").Value = 200046400
$ws.Range("J55").Value = 200046400
$ws.Range("L55").Value = 200046400
$ws.Range("N55").Value = -200047030

# Row 80
$ws.Range("H80").Value = 125025360
$ws.Range("I80").Value = 20833.334
$ws.Range("J80").Value = 200028080
$ws.Range("K80").Value = 20833.334
$ws.Range("L80").Value = 200028080
$ws.Range("M80").Value = -19835.334
$ws.Range("N80").Value = -200030076

# Row 83
$ws.Range("H83").Value = 125025360
$ws.Range("I83").Value = 20833.334
$ws.Range("J83").Value = 200028080
$ws.Range("K83").Value = 62500.00199999999
$ws.Range("L83").Value = 600084240
$ws.Range("M83").Value = -57508.00199999999
$ws.Range("N83").Value = -600094224

# Row 97
$ws.Range("H97").Value = 747.3158
$ws.Range("I97").Value = 747.3158
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 747.3158
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -251.3158
$ws.Range("N97").ClearContents()

# Row 109
$ws.Range("H109").Value = 35188
$ws.Range("I109").Value = 34999
$ws.Range("J109").Value = 35377
$ws.Range("K109").Value = 34999
$ws.Range("L109").Value = 35377
$ws.Range("M109").Value = -33612
$ws.Range("N109").Value = -38151

$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Range("H108").Value = 38684
$ws.Range("J108").Value = 38684
$ws.Range("L108").Value = 38684
$ws.Range("N108").Value = -46364

$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 2662.5
$ws.Range("I8").Value = 883.3333
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 883.3333
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -743.3333
$ws.Range("N8").Value = -8280

# Row 31
$ws.Range("H31").Value = 5183.091
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5183.091
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5183.091
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5773.091

# Row 34
$ws.Range("H34").Value = 5183.091
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5183.091
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5183.091
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5587.091

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 210.33333
$ws.Range("J17").Value = 210.33333
$ws.Range("L17").Value = 630.99999
$ws.Range("N17").Value = -968.99999

# Row 26
$ws.Range("H26").Value = 659.4545000000001
$ws.Range("I26").Value = 66.25
$ws.Range("K26").Value = 198.75
$ws.Range("M26").Value = 89.25

# Row 34
$ws.Range("H34").Value = 3360.375
$ws.Range("I34").Value = 255.14285
$ws.Range("J34").Value = 5775.5557
$ws.Range("K34").Value = 765.4285500000001
$ws.Range("L34").Value = 17326.6671
$ws.Range("M34").Value = -681.4285500000001
$ws.Range("N34").Value = -17494.6671

# Row 38
$ws.Range("H38").Value = 80.35294
$ws.Range("I38").Value = 22.5
$ws.Range("J38").Value = 111.90909
$ws.Range("K38").Value = 67.5
$ws.Range("L38").Value = 335.72727
$ws.Range("M38").Value = 279.5
$ws.Range("N38").Value = -1029.72727

# Row 55
$ws.Range("H55").Value = 4962.5
$ws.Range("J55").Value = 4962.5
$ws.Range("L55").Value = 14887.5
$ws.Range("N55").Value = -15241.5

# Row 56
$ws.Range("H56").Value = 87199.914
$ws.Range("I56").Value = 87199.914
$ws.Range("K56").Value = 87199.914
$ws.Range("M56").Value = -86669.914

# Row 131
$ws.Range("H131").Value = 16705.033
$ws.Range("I131").Value = 389.89795
$ws.Range("J131").Value = 96649.2
$ws.Range("K131").Value = 1169.69385
$ws.Range("L131").Value = 289947.6
$ws.Range("M131").Value = 3870.30615
$ws.Range("N131").Value = -300027.6

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1437.6
$ws.Range("I97").Value = 1200.409
$ws.Range("J97").Value = 3177
$ws.Range("K97").Value = 1200.409
$ws.Range("L97").Value = 3177
$ws.Range("M97").Value = -704.4090000000001
$ws.Range("N97").Value = -4169

# Row 113
$ws.Range("H113").Value = 4870.3335
$ws.Range("I113").Value = 4870.3335
$ws.Range("K113").Value = 4870.3335
$ws.Range("M113").Value = -2700.3335

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 940.9091
$ws.Range("I22").Value = 1450
$ws.Range("J22").Value = 890
$ws.Range("K22").Value = 1450
$ws.Range("L22").Value = 890
$ws.Range("M22").Value = -1155
$ws.Range("N22").Value = -1480

# Row 27
$ws.Range("H27").Value = 940.9091
$ws.Range("I27").Value = 1450
$ws.Range("J27").Value = 890
$ws.Range("K27").Value = 1450
$ws.Range("L27").Value = 890
$ws.Range("M27").Value = -1343
$ws.Range("N27").Value = -1104

# Row 46
$ws.Range("H46").Value = 856
$ws.Range("I46").Value = 670
$ws.Range("K46").Value = 670
$ws.Range("M46").Value = -482

# Row 55
$ws.Range("H55").Value = 143426.58
$ws.Range("I55").Value = 222770.11
$ws.Range("J55").Value = 608.2
$ws.Range("K55").Value = 222770.11
$ws.Range("L55").Value = 608.2
$ws.Range("M55").Value = -222597.11
$ws.Range("N55").Value = -954.2

# Row 100
$ws.Range("H100").Value = 4387.8335
$ws.Range("I100").Value = 2665.0833
$ws.Range("K100").Value = 2665.0833
$ws.Range("M100").Value = -2124.0833

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 10433.333
$ws.Range("I29").Value = 8300
$ws.Range("J29").Value = 11500
$ws.Range("K29").Value = 8300
$ws.Range("L29").Value = 11500
$ws.Range("M29").Value = -8010
$ws.Range("N29").Value = -12080

# Row 62
$ws.Range("H62").Value = 3895.1765
$ws.Range("I62").Value = 3760
$ws.Range("J62").Value = 3951.5
$ws.Range("K62").Value = 3760
$ws.Range("L62").Value = 3951.5
$ws.Range("M62").Value = -3136
$ws.Range("N62").Value = -5199.5

# Row 65
$ws.Range("H65").Value = 3895.1765
$ws.Range("I65").Value = 3760
$ws.Range("J65").Value = 3951.5
$ws.Range("K65").Value = 18800
$ws.Range("L65").Value = 19757.5
$ws.Range("M65").Value = -15680
$ws.Range("N65").Value = -25997.5

# Row 70
$ws.Range("H70").Value = 31701.666
$ws.Range("J70").Value = 31701.666
$ws.Range("L70").Value = 31701.666
$ws.Range("N70").Value = -32331.666

# Row 73
$ws.Range("H73").Value = 31701.666
$ws.Range("J73").Value = 31701.666
$ws.Range("L73").Value = 31701.666
$ws.Range("N73").Value = -33885.666

# Row 94
$ws.Range("H94").Value = 124000
$ws.Range("J94").Value = 124000
$ws.Range("L94").Value = 124000
$ws.Range("N94").Value = -125802

# Row 96
$ws.Range("H96").Value = 3101
$ws.Range("I96").Value = 3101
$ws.Range("K96").Value = 3101
$ws.Range("M96").Value = -1728

# Row 101
$ws.Range("H101").Value = 4980
$ws.Range("J101").Value = 4980
$ws.Range("L101").Value = 4980
$ws.Range("N101").Value = -11470

# Row 103
$ws.Range("H103").Value = 42840.8
$ws.Range("J103").Value = 42840.8
$ws.Range("L103").Value = 42840.8
$ws.Range("N103").Value = -45184.8

# Row 104
$ws.Range("H104").Value = 39913.332
$ws.Range("J104").Value = 39913.332
$ws.Range("L104").Value = 39913.332
$ws.Range("N104").Value = -46901.332

# Row 106
$ws.Range("H106").Value = 25188.5
$ws.Range("J106").Value = 25188.5
$ws.Range("L106").Value = 25188.5
$ws.Range("N106").Value = -27712.5

# Row 107
$ws.Range("H107").Value = 2413.5
$ws.Range("I107").Value = 651
$ws.Range("J107").Value = 3001
$ws.Range("K107").Value = 1953
$ws.Range("L107").Value = 9003
$ws.Range("M107").Value = -33
$ws.Range("N107").Value = -12843

# Row 109
$ws.Range("H109").Value = 59133.332
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59133.332
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59133.332
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -61907.332

# Row 112
$ws.Range("H112").Value = 54663.332
$ws.Range("J112").Value = 54663.332
$ws.Range("L112").Value = 54663.332
$ws.Range("N112").Value = -57617.332

# Row 113
$ws.Range("H113").Value = 1285.5333
$ws.Range("I113").Value = 853.7778
$ws.Range("J113").Value = 1933.1666
$ws.Range("K113").Value = 2561.3334
$ws.Range("L113").Value = 5799.4998
$ws.Range("M113").Value = -391.3334
$ws.Range("N113").Value = -10139.4998
